$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 261, shifting existing rows 261..354 down to 262..355.
$ws.Rows.Item(261).Insert()

# Populate the newly inserted row 261 with the new data record.
$ws.Range("A261").Value = 11
$ws.Range("B261").Value = "Vega Monumental Concepción"
$ws.Range("C261").Value = "Bíobío"
$ws.Range("D261").Value = 45202
$ws.Range("E261").Value = 8
$ws.Range("F261").Value = 100112040
$ws.Range("G261").Value = "Cilantro"
$ws.Range("H261").Value = "Sin especificar"
$ws.Range("I261").Value = "Primera"
$ws.Range("J261").Value = 50
$ws.Range("K261").Value = 5500
$ws.Range("L261").Value = 6000
$ws.Range("M261").Value = 5800
$ws.Range("N261").Value = "$/caja 36 atados"
$ws.Range("O261").Value = "Región Metropolitana"
$ws.Range("P261").Value = 161
$ws.Range("Q261").Value = 36
$ws.Range("R261").Value = "Hortaliza"
